# Auto-generated edit script: updates crypto price/volume data per commit
# "Updated cryptos list on Sat Apr 27 10:35:55 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.994.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.124.06'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.52%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.28'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.18'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.93%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.116.65'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.63%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.46%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.85%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.76%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.72%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.06'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.634.35'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.038.99'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.127.29'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '470.92'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.08'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.42%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.17%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.65'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.65%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.90'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.71%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.84'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.63%  '

$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.08'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.00%  '

$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.82'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.98%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.08%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.108'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.18%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.16%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.41%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.26%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.91%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -10.58%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0386'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.93%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '416.94'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -6.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.20'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.897.92'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.66'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -12.19%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.262'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.50%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.09'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.38'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.45%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -7.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.08'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.15%  '

